# Generate Report for Handback
# This script fills in the "latest handback" columns (I/J/K/P) for row 7
# (the aab80b49-ac29-4516-aeee-e3326de19dd8 file) on both locale sheets,
# since a new (now-stale) handback was received for that file.

$wb = $excel.ActiveWorkbook

$latestHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0670a015e5ee8d7c97ffbb64ac51025908834074/e2e/aab80b49-ac29-4516-aeee-e3326de19dd8.md"
$latestHandbackDisplay = "aab80b49-ac29-4516-aeee-e3326de19dd8.md"
$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3ea62654b411e73b5cfc446484869cda254c7ba/e2e/aab80b49-ac29-4516-aeee-e3326de19dd8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0670a015e5ee8d7c97ffbb64ac51025908834074/e2e/aab80b49-ac29-4516-aeee-e3326de19dd8.md."

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = $latestHandbackDisplay
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestHandbackUrl, "", "", $latestHandbackDisplay) | Out-Null
$wsZh.Range("I7").Style = "HyperLink"

$wsZh.Range("J7").Value = "aab80b49-ac29-4516-aeee-e3326de19dd8.9a5703b60ea8578c1b069a74d9273734c74b56df.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-15 22:51:57"
$wsZh.Range("P7").Value = $errorMessage

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = $latestHandbackDisplay
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestHandbackUrl, "", "", $latestHandbackDisplay) | Out-Null
$wsDe.Range("I7").Style = "HyperLink"

$wsDe.Range("J7").Value = "aab80b49-ac29-4516-aeee-e3326de19dd8.9a5703b60ea8578c1b069a74d9273734c74b56df.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-15 22:52:10"
$wsDe.Range("P7").Value = $errorMessage

Write-Output "Report generated for handback."
